# Refitting NCDEs to individual patients (for manuscript figure)
# Adds a new "Label" column (H) to the classification results sheet and
# updates some refit prediction/error values (columns D and E, plus F11).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Add new "Label" header in H1, matching the style of the other headers ---
$ws.Range("G1").Copy()
$ws.Range("H1").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = 0
$ws.Range("H1").Value = "Label"

# --- Block 1 (Iterations = 100), rows 2-11 ---
$ws.Range("D2").Value = 0.2582680359074394
$ws.Range("E2").Value = 0.2582680359074394
$ws.Range("H2").Value = 0

$ws.Range("D3").Value = 0.6847620740676935
$ws.Range("E3").Value = 0.6847620740676935
$ws.Range("H3").Value = 0

$ws.Range("H4").Value = 0

$ws.Range("D5").Value = 0.5980874736696493
$ws.Range("E5").Value = 0.5980874736696493
$ws.Range("H5").Value = 0

$ws.Range("D6").Value = 0.6925896487871037
$ws.Range("E6").Value = 0.6925896487871037
$ws.Range("H6").Value = 0

$ws.Range("D7").Value = 0.5979440315102105
$ws.Range("E7").Value = 0.4020559684897895
$ws.Range("H7").Value = 1

$ws.Range("D8").Value = 0.604387125276084
$ws.Range("E8").Value = 0.395612874723916
$ws.Range("H8").Value = 1

$ws.Range("D9").Value = 0.3718158402351632
$ws.Range("E9").Value = 0.6281841597648368
$ws.Range("H9").Value = 1

$ws.Range("H10").Value = 1

$ws.Range("D11").Value = 0.7883496883457821
$ws.Range("E11").Value = 0.2116503116542179
$ws.Range("F11").Value = 0.7012858390808105
$ws.Range("H11").Value = 1

# --- Block 2 (Iterations = 200), rows 12-21 : only the new Label column ---
$ws.Range("H12").Value = 0
$ws.Range("H13").Value = 0
$ws.Range("H14").Value = 0
$ws.Range("H15").Value = 0
$ws.Range("H16").Value = 0
$ws.Range("H17").Value = 1
$ws.Range("H18").Value = 1
$ws.Range("H19").Value = 1
$ws.Range("H20").Value = 1
$ws.Range("H21").Value = 1
